$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet originally had a two-row (pandas multi-index) header: row 1 held the
# top-level group labels ("Tackles", "Challenges", "Blocks", "Unnamed: N_level_0"
# placeholders) and row 2 held the real column names. Clean this up by promoting
# real, unambiguous column headers into row 1, then hiding the old helper rows
# (2, 3, and the summary row 20) instead of deleting them, and removing the
# now-unnecessary header merges.

# Remove the column-group merges across the old row-1 header.
$ws.Range("H1:L1").UnMerge()
$ws.Range("M1:P1").UnMerge()
$ws.Range("Q1:S1").UnMerge()

# Write the cleaned-up header labels into row 1.
$ws.Range("A1").Value = "Player ID"
$ws.Range("B1").Value = "Player"
$ws.Range("C1").Value = "#"
$ws.Range("D1").Value = "Nation"
$ws.Range("E1").Value = "Pos"
$ws.Range("F1").Value = "Age"
$ws.Range("G1").Value = "90s"
$ws.Range("H1").Value = "Tkl"
$ws.Range("I1").Value = "TklW"
$ws.Range("J1").Value = "Def 3rd"
$ws.Range("K1").Value = "Mid 3rd"
$ws.Range("L1").Value = "Att 3rd"
$ws.Range("M1").Value = "Cha"
$ws.Range("N1").Value = "Att"
$ws.Range("O1").Value = "Tkl%"
$ws.Range("P1").Value = "Lost"
$ws.Range("Q1").Value = "Blocks"
$ws.Range("R1").Value = "Sh"
$ws.Range("S1").Value = "Pass"
$ws.Range("T1").Value = "Int"
$ws.Range("U1").Value = "Tkl+Int"
$ws.Range("V1").Value = "Clr"
$ws.Range("W1").Value = "Err"

# Fill in the missing Tkl% values (players with 0 tackle attempts had a blank
# cell instead of a 0).
$ws.Range("O7").Value = 0
$ws.Range("O14").Value = 0
$ws.Range("O18").Value = 0

# Hide the now-redundant old header row, the blank spacer row, and the
# aggregate "16 Players" summary row instead of removing them.
$ws.Rows.Item(2).Hidden = $true
$ws.Rows.Item(3).Hidden = $true
$ws.Rows.Item(20).Hidden = $true
